# Loan RBI, Variable Instalments
# Inserts a new (blank) column before column N ("Late") on the
# "Repayment Schedule" sheet, widens it to match the other data columns,
# and moves the active sheet/selection to that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before N (shifts old N/O/P -> O/P/Q).
$ws.Columns("N:N").Insert()

# Match the column width used by the diff (width="10", customWidth only).
$ws.Columns("N:N").ColumnWidth = 9.1

# Make "Repayment Schedule" the active sheet/tab and move the selection.
$ws.Activate()
$ws.Range("S10").Select()
